# Screenshots Required To Do List - add "Weapons In Inventory" and
# "Sort Weapons" screenshot rows just above the existing "IHasSummary
# interface" row (i.e. insert two new rows at row 12, pushing the
# remaining rows down by two).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at 12:13 - this shifts rows 12-25 down to 14-27
# and re-bases the sheet dimension automatically.
$ws.Rows("12:13").Insert()

# The inserted rows come back with no explicit cell-level style (only
# the column default), whereas every other data row in this sheet
# carries an explicit style index on column A (s=4) and column B
# (s=3). Copy the formatting down from the row directly above (row 11)
# so the new rows match the look of the rest of the table before we
# fill in their text.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)

# Row 12: weapon inventory screenshot requirement
$ws.Range("B12").Value = "User having the functionality to sort the damage of the weapons"
$ws.Range("B13").Value = "User Sorted the weapons by descending damage"

# Filenames for the two new screenshots
$ws.Range("C12").Value = "weapons in inventory.png"
$ws.Range("C13").Value = "weapons sort by descending damage.png"

# Match the author's final selection position.
$ws.Range("D24").Select()
